$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell "Save" in H1, matching the style/format of existing headers (e.g. G1)
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null
$ws.Range("H1").Value = "Save"

# Add new data cell H2 = 0 (plain numeric, no special style)
$ws.Range("H2").Value = 0
